$d = $word.ActiveDocument

# The title paragraph originally reads (run-by-run):
#   "{{ caseType }}" <br> "{{ caseId }" | "}{% if " | "dossierNr" | " %} ({{ "
#   | "dossierNr" | " }}){% endif %}" | <br> | "{{ formType }} " <br>
# (the "|" marks a run boundary between runs that share identical rPr).
#
# 1) Join "{{ caseId }" and "}{% if " (split across two runs with identical
#    formatting) back into one run / one piece of text: "{{ caseId }}{% if ".
#    A Find/Replace across the run boundary merges the matched runs because
#    they carry the same character formatting.
$d.Content.Find.Execute("{{ caseId }" + "}{% if ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ caseId }}{% if ", 2)

# 2) Drop the run break between " }}){% endif %}" and the manual line break
#    that follows it (that run only ever contained the <w:br/>). Re-matching
#    the text together with the line-break character (Chr(11)) and replacing
#    it with itself causes the two adjacent, identically-formatted runs to
#    be merged into a single run again.
$lineBreak = [char]11
$d.Content.Find.Execute(" }}){% endif %}" + $lineBreak, $true, $false, $false, $false, $false, `
                         $true, 1, $false, " }}){% endif %}" + $lineBreak, 2)

# 3) Simplify the trailing "{{ formType }} " text down to a single space.
$d.Content.Find.Execute("{{ formType }} ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " ", 2)
